# "Reorganización completa: limpieza de módulos antiguos, nuevas entregas y optimización"
#
# 1) Rename the sheet "Datos" -> "produccion_leche" (in place, keeps sheetId/r:id).
# 2) Strip the custom header styling (bold white font on a solid blue fill,
#    centered alignment) from A1:F1, reverting the header row to the
#    workbook's default (unstyled) cell format.
# 3) Reset the custom 20-character column widths on A:F back to Excel's
#    standard/default column width.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "produccion_leche"

$headerRange = $ws.Range("A1:F1")
$headerRange.ClearFormats()

$dataCols = $ws.Columns("A:F")
$dataCols.ColumnWidth = $ws.StandardWidth
